# version 1.0.3: fixed bug introduced with v1.0.1 in logistical data extraction
#
# The previous version (v1.0.1) had started writing extracted logistics
# data (shipper plantcode, destination material/quantity rows) into cells
# that should remain empty user-input placeholders. This removes that
# erroneously-written data and restores the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Remove the erroneously populated values.
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("A10:B10").ClearContents()
$ws.Range("A11:B11").ClearContents()

# Restore the active selection/cursor position recorded at save time.
$ws.Range("E11").Select()
